# Weekly fruit/vegetable price update: a new weekly observation was
# inserted into the "Haba" (broad bean) price series at row 379, pushing
# all subsequent rows (old 379-410) down by one (new 380-411).
#
# This mirrors how the source data pipeline appends/prepends a new
# week's record into the middle of the date-sorted series.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 379 (shifts 379..410 down to 380..411,
# and Excel copies formatting such as the date style from the row above).
$ws.Rows.Item(379).Insert()

# Populate the newly inserted row with the new observation.
$ws.Range("A379").Value = 9
$ws.Range("B379").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C379").Value = "Metropolitana"
$ws.Range("D379").Value = 45223
$ws.Range("E379").Value = 13
$ws.Range("F379").Value = 100112026
$ws.Range("G379").Value = "Haba"
$ws.Range("H379").Value = "Sin especificar"
$ws.Range("I379").Value = "Primera"
$ws.Range("J379").Value = 70
$ws.Range("K379").Value = 10000
$ws.Range("L379").Value = 12000
$ws.Range("M379").Value = 11000
$ws.Range("N379").Value = "$/saco 25 kilos"
$ws.Range("O379").Value = "Provincia de Melipilla"
$ws.Range("P379").Value = 440
$ws.Range("Q379").Value = 25
$ws.Range("R379").Value = "Hortaliza"
